$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of real-time GDP data to insert before the existing row 2
# (backward extension of the series), shifting all existing rows down by 11.
$newData = @(
    @(30681, 1.466797881812631),
    @(31047, 2.900424903011278),
    @(31412, 2.60323159784559),
    @(31777, 2.279090113735815),
    @(32142, 1.278816132757399),
    @(32508, 3.441722972972983),
    @(32873, 4.033476219636634),
    @(33238, 5.482086096613425),
    @(33603, 5.2529761904762),
    @(33969, 1.60115933832885),
    @(34334, -0.9914767785701772)
)

$n = $newData.Count

# Insert n new blank rows above the current row 2, pushing existing data down.
$insertRange = $ws.Range("A2:B$($n + 1)")
$insertRange.EntireRow.Insert()

# The inserted rows picked up a blended style from row 1 on insert; reset to
# match the rest of the table: column A uses the date-format style, column B
# is unstyled (same as every other data row).
$ws.Range("B2:B$($n + 1)").ClearFormats()
$ws.Range("A$($n + 2)").Copy()
$ws.Range("A2:A$($n + 1)").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the newly inserted rows with the extended historical data.
for ($i = 0; $i -lt $n; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 1).Value = $newData[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $newData[$i][1]
}
